# Upload new version with timestamp
# Populates the first data row of the report with a transaction line and
# marks the two "text-ish" columns (item name / current balance) as
# Text-formatted (numFmtId 49, i.e. "@") so values like "0:0" / "1:0"
# are preserved verbatim instead of being re-interpreted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: new data line ----------------------------------------------
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "CEFOTAX 2 GM VIAL"   # merged B4:G4 - item name
$ws.Range("H4").Value = "0:0"                  # merged H4:K4 - current balance
$ws.Range("L4").Value = 62                     # merged L4:M4 - selling price
$ws.Range("N4").Value = "1:0"                  # number of transactions

# --- Row 5: totals / balance carry row ----------------------------------
$ws.Range("K5").Value = 62                     # merged K5:N5

# Row 5 grew slightly taller to fit the new value.
$ws.Rows(5).RowHeight = 26.25

# --- Number formats ------------------------------------------------------
# The "name" column (B4:G4 merged + N4 share style s=7) and the "balance"
# column (H4:K4 merged, style s=8) both need to be stored as Text so
# values such as "0:0" / "1:0" round-trip literally.
$ws.Range("B4:G4").NumberFormat = "@"
$ws.Range("N4").NumberFormat = "@"
$ws.Range("H4:K4").NumberFormat = "@"
